# Refresh the cryptos list (Price + Volume(1h) columns) with the latest
# scraped figures. Price cells that happen to look like plain numbers are
# forced to text (NumberFormat "@") before the assignment, then reset to
# "General" formatting, so they stay literal strings like the source data
# (e.g. "0.999", "314.61") instead of being auto-converted to numeric
# cells by Excel's type inference. Values with two dots (e.g. "46.287.65")
# are never number-like, so they don't need this treatment.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.287.65"
$ws.Range("E2").Value = "  +1.32%  "

$ws.Range("D3").Value = "2.613.85"
$ws.Range("E3").Value = "  +7.98%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.61"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +4.60%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.80"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +3.65%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.600"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +5.83%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("E9").Value = "  +12.86%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.80"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +11.53%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.32"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +0.36%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0841"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  +5.87%  "

$ws.Range("E13").Value = "  +15.66%  "

$ws.Range("D14").Value = "3.009.13"
$ws.Range("E14").Value = "  +7.94%  "

$ws.Range("E15").Value = "  +1.46%  "

$ws.Range("D16").Value = "2.608.30"
$ws.Range("E16").Value = "  +7.21%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.913"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  +7.78%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "15.18"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +5.36%  "

$ws.Range("D19").Value = "46.481.32"
$ws.Range("E19").Value = "  +1.76%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.35"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +0.98%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000102"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +7.26%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.75"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +8.52%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.88"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +4.82%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "254.99"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +4.49%  "

$ws.Range("E25").Value = "  +9.68%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.21"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  +13.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.96"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +30.22%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +0.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.64"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  +8.33%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "40.74"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +4.65%  "

$ws.Range("E31").Value = "  +2.13%  "

$ws.Range("E32").Value = "  +10.55%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.74"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -2.41%  "

$ws.Range("E34").Value = "  +13.84%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.87"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +4.63%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "153.92"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +3.75%  "

$ws.Range("E37").Value = "  +7.68%  "

$ws.Range("E38").Value = "  +4.76%  "

$ws.Range("E39").Value = "  +5.86%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.18"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +12.55%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.24"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +8.41%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.64"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  +10.46%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0328"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +8.48%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.76"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +33.23%  "

$ws.Range("D45").Value = "2.035.04"
$ws.Range("E45").Value = "  +4.29%  "

$ws.Range("E46").Value = "  -0.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "91.12"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -0.06%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "112.52"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +10.45%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.24"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +6.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.79"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +1.81%  "

$ws.Range("E51").Value = "  +8.12%  "
